$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "O" (Oxygen) current value in C2 with the new measurement.
# Dependent formulas in D2 and F2 will recalculate automatically.
$ws.Range("C2").Value = [double]"1.1657810000000001E+30"

# Move the active selection to D2 to match the saved selection state.
$ws.Activate()
$ws.Range("D2").Select()
